# Change font sizes for Aspuru-Guzik figure
# Updates the "k = 2 partition" / "k = 3 partition" labels:
#  - reposition/resize the text boxes
#  - center-align the paragraph
#  - shrink the font from 36pt to 30pt
#  - pluralize "partition" -> "partitions"

$emuPerPt = 12700

# The COM position/size properties round-trip EMU through a
# single-precision (float32) point value; converting back to EMU on
# save can truncate and land 1 EMU short of the intended value. Nudge
# by a tiny epsilon (well under 1/12700 pt) so the truncated result
# still lands on the exact EMU we want.
function EmuToPt($emu) {
    return ($emu / $emuPerPt) + 0.00005
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$k2 = $s.Shapes.Item("TextBox 65")
$k3 = $s.Shapes.Item("TextBox 66")

# --- "k = 2 partition" textbox ---
$k2.TextFrame.TextRange.Text = "k = 2 partitions"
$k2.TextFrame.TextRange.Font.Size = 30
$k2.TextFrame.TextRange.ParagraphFormat.Alignment = 2  # ppAlignCenter

$k2.Left = EmuToPt 1761087
$k2.Top = EmuToPt 1027582
$k2.Width = EmuToPt 2517036
$k2.Height = EmuToPt 553998

# --- "k = 3 partition" textbox ---
$k3.TextFrame.TextRange.Text = "k = 3 partitions"
$k3.TextFrame.TextRange.Font.Size = 30
$k3.TextFrame.TextRange.ParagraphFormat.Alignment = 2  # ppAlignCenter

$k3.Left = EmuToPt 6182213
$k3.Top = EmuToPt 1027582
$k3.Width = EmuToPt 2517036
$k3.Height = EmuToPt 553998
